$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 (existing rows 7-31 shift down to 8-32)
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new weekly price entry
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(7, 3).Value = "Metropolitana"
$ws.Cells.Item(7, 4).Value = 44859
$ws.Cells.Item(7, 5).Value = 13
$ws.Cells.Item(7, 6).Value = 300000001
$ws.Cells.Item(7, 7).Value = "Rabanito"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 7900
$ws.Cells.Item(7, 11).Value = 3000
$ws.Cells.Item(7, 12).Value = 3000
$ws.Cells.Item(7, 13).Value = 3000
$ws.Cells.Item(7, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(7, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(7, 16).Value = 30
$ws.Cells.Item(7, 17).Value = 100
$ws.Cells.Item(7, 18).Value = "Hortaliza"
